$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "27.031.96"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "1.651.98"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "215.14"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "19.87"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "1.886.22"
$ws.Range("D13").Value = "1.654.35"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").Value = "65.33"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "240.19"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "27.017.06"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +4.47%  "
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("E24").Value = "  +3.54%  "
$ws.Range("D25").Value = "145.74"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").Value = "15.82"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").Value = "1.521.22"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("E35").Value = "  +8.67%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "0.581"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  +8.66%  "
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("D43").Value = "65.84"
$ws.Range("E43").Value = "  +8.87%  "
$ws.Range("D44").Value = "1.792.40"
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "89.50"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "0.0976"
$ws.Range("E51").Value = "  +2.06%  "
